$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Wine & naming convention" sheet: append the PDO/PGI equivalence table
#    and the "Millesime" row (rows 10-19). Cells are written in the exact
#    order the original authors entered them so new shared strings land at
#    the same indices as the target workbook.
# ---------------------------------------------------------------------------
$wsNaming = $wb.Worksheets.Item("Wine & naming convention")

$wsNaming.Range("A10").Value = "Bottle's label mandatory information"
$wsNaming.Range("B10").Value = "Denomination; Origin; Alcohol percentage; Bottler identity (or importer); Volume; Lot number; Pregnancy avertissement; Allergen information"

$wsNaming.Range("A11").Value = "PDO equivalent in France"
$wsNaming.Range("B11").Value = "AOC - Appellation d’origine contrôlée"

$wsNaming.Range("B12").Value = "VDQS - Vin délimité de qualité supérieure."

$wsNaming.Range("A13").Value = "PDO equivalent in Italy"

$wsNaming.Range("A12").Value = "PGI equivalent in France"

$wsNaming.Range("A14").Value = "PGI equivalent in Italy"

$wsNaming.Range("B13").Value = "DOC -  Denominazione di Origine Controllata"
$wsNaming.Range("B14").Value = "DOCG - Denominazione di Origine Controllata e Garantita"

$wsNaming.Range("A15").Value = "PDO equivalent in Spain"
$wsNaming.Range("A16").Value = "PGI equivalent in Spain"

$wsNaming.Range("B15").Value = "DO - Denominaciòn de Origen"
$wsNaming.Range("B16").Value = "DOCa - Denominaciòn d’Origen Calificada"

$wsNaming.Range("A17").Value = "PDO equivalent in Portugal"
$wsNaming.Range("A18").Value = "PDO equivalent in Germany"

$wsNaming.Range("B17").Value = "DOC - Denominação de Origem Controlada"
$wsNaming.Range("B18").Value = "QBA - Qualitätswein bestimmter Anbaugebiete"

$wsNaming.Range("A19").Value = "Millésime"
$wsNaming.Range("B19").Value = "A wine can have a millesime mention if 95% of the grapes come from the same year in USA, and 85% in Europe"

# The "Millesime" label wraps, like the other short emphasised labels on this sheet.
$wsNaming.Range("A19").WrapText = $true

$wsNaming.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$wsNaming.Range("A20").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) "Wine tasting" sheet: append "Tasting order" at the end of the table,
#    then insert the "Spitting" row in the middle (row 25), pushing the
#    existing rows 25-34 down to 26-35.
# ---------------------------------------------------------------------------
$wsTasting = $wb.Worksheets.Item("Wine tasting")

$wsTasting.Range("A34").Value = "Tasting order"
$wsTasting.Range("B34").Value = "If you taste different wine, here is the recommended order: 1. Sparkling wine; 2. White wine, rosés; 3. Young red; 4. Strong white or sweet wine; 5. Strong red wine; 6. Fortified wine"

$wsTasting.Rows(25).Insert()
$wsTasting.Range("A25").Value = "Spitting"
$wsTasting.Range("B25").Value = "Spitting is recommended to stay fully lucid while tasting different wine"

# This becomes the active sheet/tab of the workbook.
$wsTasting.Activate() | Out-Null
$wsTasting.Range("B6").Select() | Out-Null
